$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the "statut_label" value: "bleu" -> "noir"
$ws.Cells.Replace("bleu", "noir", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)

# Correct the "statut_name" values (fix wording)
$ws.Cells.Replace("pas de résultat ni de publication", "pas de résultat postés ni publiés", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$ws.Cells.Replace("résultat et / ou publication posté dans les 36 mois", "résultat postés ou publiés dans les 36 mois", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$ws.Cells.Replace("résultat et / ou publication posté dans les 12 mois", "résultat postés ou publiés dans les 12 mois", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
$ws.Cells.Replace("résultat et / ou publication posté", "résultat postés ou publiés", [Microsoft.Office.Interop.Excel.XlLookAt]::xlWhole)
